$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 holds the rule-name label for the last rule row (was "R40"). The
# column stores text (shared-string) values, so force text formatting
# before writing "1" - otherwise Excel's usual smart-typing would store
# it as the number 1 instead of the string "1".
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
